# Apply "Add same ttime partners" edit to watch_list workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$NL = [char]10

# --- Column E width: 15 -> 16 (match column D / F width) ---
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# --- Row 44 ---
$ws.Range("D44").Value = "אור" + $NL + "דורון"
$ws.Range("E44").Value = "קריספין" + $NL + "רווה"

# --- Row 46 ---
$ws.Range("F46").Value = "דבוש" + $NL + "פיאצה"

# --- Row 47 ---
$ws.Range("C47").Value = "אסף" + $NL + "אסרף"

# --- Row 50 ---
$ws.Range("C50").Value = "לוטם" + $NL + "סיני"
$ws.Range("D50").Value = "אבנר" + $NL + "לומיאנסקי"
$ws.Range("E50").Value = "דעאל" + $NL + "שגיא"
$ws.Range("F50").Value = "איתי כהן" + $NL + "כלפה"

# --- Row 53 ---
$ws.Range("C53").Value = "אנזו" + $NL + "לואיס"
$ws.Range("D53").Value = "אנדי" + $NL + "דוד"

# --- Row 56 ---
$ws.Range("C56").Value = "ארד" + $NL + "יואל"
$ws.Range("E56").Value = "דורון" + $NL + "שמעון"

# --- Row 59 ---
$ws.Range("C59").Value = "קריספין" + $NL + "רווה"
$ws.Range("D59").Value = "דבוש" + $NL + "פיאצה"
$ws.Range("E59").Value = "אור" + $NL + "מרדש"

# --- Row 62 ---
$ws.Range("C62").Value = "שראל" + $NL + "שרעבי"
$ws.Range("D62").Value = "יונג" + $NL + "ניסנוב"
$ws.Range("E62").Value = "דימנטמן" + $NL + "מטמוני"

# --- Row 65 ---
$ws.Range("D65").Value = "דותן" + $NL + "שגיא"
$ws.Range("E65").Value = "אסף" + $NL + "אסרף"
